$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '44.880.37'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.85%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.257.82'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.86%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.79%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '300.81'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.57'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.565'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.88%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.59%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.507'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.02'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.29%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0785'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.14%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.16'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.54%  '
$ws.Range("E13").Value = '  -0.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.601.52'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.260.91'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.64'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.64%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.796'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -4.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '44.703.03'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.62'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +6.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0915'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.89%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.07'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.06'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.51%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '238.79'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.85%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.87'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.84%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.89'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -3.83%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.28'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.26%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '39.30'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +4.75%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.55'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.79%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.54'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.83%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '153.05'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.53%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.57'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -6.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0776'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.55'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.65%  '
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.92'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.30%  '
$ws.Range("B36").Value = 'Stellar'
$ws.Range("C36").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.117'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.48%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.106'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -4.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.72'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -6.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.25'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.70'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.82'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -8.46%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("D43").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.790.90'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.41%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.82'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +8.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.188'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '69.88'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.51%  '
$ws.Range("B48").Value = 'BitcoinSV'
$ws.Range("C48").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '75.93'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.40%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '96.02'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.67'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -4.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.81'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.63%  '
